$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 198, shifting rows 198:202 down to 199:203
$ws.Rows.Item(198).Insert()

# Populate the new row 198 with the new record's data
$ws.Cells.Item(198, 1).Value = 3
$ws.Cells.Item(198, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(198, 3).Value = "Coquimbo"
$ws.Cells.Item(198, 4).Value = 44448
$ws.Cells.Item(198, 5).Value = 5
$ws.Cells.Item(198, 6).Value = 100112032
$ws.Cells.Item(198, 7).Value = "Zapallo italiano"
$ws.Cells.Item(198, 8).Value = "Sin especificar"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 160
$ws.Cells.Item(198, 11).Value = 13000
$ws.Cells.Item(198, 12).Value = 14000
$ws.Cells.Item(198, 13).Value = 13516
$ws.Cells.Item(198, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(198, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(198, 16).Value = 193
$ws.Cells.Item(198, 17).Value = 70
$ws.Cells.Item(198, 18).Value = "Hortaliza"
